$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Object_Mapping")

# B3: Electrolyzer_PEM -> PEM_Electrolyzer
$ws.Range("B3").Value = "PEM_Electrolyzer"

# B5: Methanol_Plant -> Destilation_tower
$ws.Range("B5").Value = "Destilation_tower"
